$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N (14th column),
# shifting the "Late" and "Outstanding" columns one to the right.
$ws.Columns("N").Insert()

# Match Excel's default behavior of copying the width from the column
# to the left (M) onto the freshly inserted column.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Update the active selection to match the saved view state.
$ws.Range("R8").Select()
